$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: AO2020Mar values updated
$ws.Range("A2").Value = "AO2020Mar"
$ws.Range("B2").Value = -3032.4142
$ws.Range("C2").Value = 1572.515
$ws.Range("D2").Value = -1.928
$ws.Range("E2").Value = 0.054
$ws.Range("F2").Value = -6114.487
$ws.Range("G2").Value = 49.658

# Row 3: now ma.S.L12 (was LS2020Apr)
$ws.Range("A3").Value = "ma.S.L12"
$ws.Range("B3").Value = -0.6571
$ws.Range("C3").Value = 0.08400000000000001
$ws.Range("D3").Value = -7.84
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = -0.821
$ws.Range("G3").Value = -0.493

# Row 4: now sigma2 (was ma.S.L12)
$ws.Range("A4").Value = "sigma2"
$ws.Range("B4").Value = 37290000
$ws.Range("C4").Value = 0.102
$ws.Range("D4").Value = 365000000
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 37300000
$ws.Range("G4").Value = 37300000

# Row 5 (old sigma2 row) is removed entirely
$ws.Range("A5:G5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
